# Auto-generated script to update cryptos.xlsx price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.178.77'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.234.02'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '295.02'
$ws.Range("E5").Value = '  +1.41%  '
$ws.Range("D6").Value = '88.02'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("D7").Value = '0.513'
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -0.33%  '
$ws.Range("D10").Value = '30.70'
$ws.Range("E10").Value = '  -0.42%  '
$ws.Range("D11").Value = '50.97'
$ws.Range("E11").Value = '  +6.89%  '
$ws.Range("D12").Value = '0.0784'
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("E13").Value = '  +2.76%  '
$ws.Range("D14").Value = '6.48'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = '2.580.78'
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("D16").Value = '13.88'
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").Value = '2.272.73'
$ws.Range("E17").Value = '  +2.24%  '
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("D19").Value = '40.106.23'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").Value = '0.0₃0888'
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").Value = '11.31'
$ws.Range("E21").Value = '  -3.89%  '
$ws.Range("D22").Value = '5.80'
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").Value = '65.89'
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").Value = '237.23'
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").Value = '2.49'
$ws.Range("E26").Value = '  +0.46%  '
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("D28").Value = '23.39'
$ws.Range("E28").Value = '  +3.51%  '
$ws.Range("E29").Value = '  -1.76%  '
$ws.Range("D30").Value = '9.33'
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").Value = '158.01'
$ws.Range("E31").Value = '  +3.08%  '
$ws.Range("D32").Value = '31.90'
$ws.Range("E32").Value = '  -1.09%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("E34").Value = '  +0.71%  '
$ws.Range("D35").Value = '3.07'
$ws.Range("E35").Value = '  +7.58%  '
$ws.Range("D36").Value = '0.0718'
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("E38").Value = '  +1.09%  '
$ws.Range("E39").Value = '  +3.46%  '
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("D41").Value = '15.61'
$ws.Range("E41").Value = '  -3.27%  '
$ws.Range("D42").Value = '2.099.35'
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("D43").Value = '3.72'
$ws.Range("E43").Value = '  -2.94%  '
$ws.Range("E44").Value = '  +5.59%  '
$ws.Range("D45").Value = '10.15'
$ws.Range("E45").Value = '  +2.24%  '
$ws.Range("E46").Value = '  +0.60%  '
$ws.Range("E47").Value = '  +2.58%  '
$ws.Range("E48").Value = '  -10.88%  '
$ws.Range("D49").Value = '2.449.69'
$ws.Range("E49").Value = '  +0.72%  '
$ws.Range("E50").Value = '  +3.10%  '
$ws.Range("E51").Value = '  +3.54%  '
